# Repull data / push all data: update the dSF column (F) values on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F3"  = -1
    "F4"  = -1
    "F5"  = -1
    "F6"  = 1
    "F8"  = 3
    "F9"  = 5
    "F11" = 1
    "F12" = -4
    "F13" = -1
    "F15" = -3
    "F16" = 1
    "F17" = -2
    "F18" = -1
    "F19" = 6
    "F20" = -2
    "F22" = -3
    "F23" = -1
    "F24" = 1
    "F25" = -1
    "F26" = -1
    "F27" = -1
    "F28" = -2
    "F29" = -3
    "F30" = 5
    "F31" = 0
    "F32" = 1
    "F33" = 9
    "F34" = 1
    "F35" = 1
    "F37" = -1
    "F38" = -5
    "F39" = 0
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
